$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "72.419.92"
Set-TextValue "E2" "  -0.40%  "
Set-TextValue "D3" "2.660.80"
Set-TextValue "E3" "  +0.89%  "
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "D5" "596.57"
Set-TextValue "E5" "  -1.51%  "
Set-TextValue "D6" "175.27"
Set-TextValue "E6" "  -2.23%  "
Set-TextValue "E7" "  +0.02%  "
Set-TextValue "E8" "  -0.81%  "
Set-TextValue "D9" "2.661.68"
Set-TextValue "E9" "  +1.02%  "
Set-TextValue "D10" "0.168"
Set-TextValue "E10" "  -4.30%  "
Set-TextValue "E11" "  +1.84%  "
Set-TextValue "D12" "0.355"
Set-TextValue "E12" "  -0.02%  "
Set-TextValue "D13" "5.00"
Set-TextValue "E13" "  -0.59%  "
Set-TextValue "D14" "3.149.04"
Set-TextValue "E14" "  +1.32%  "
Set-TextValue "D15" "72.390.32"
Set-TextValue "E15" "  -0.32%  "
Set-TextValue "D16" "0.0000184"
Set-TextValue "E16" "  -3.34%  "
Set-TextValue "D17" "26.21"
Set-TextValue "E17" "  -2.33%  "
Set-TextValue "D18" "2.651.16"
Set-TextValue "E18" "  +1.08%  "
Set-TextValue "D19" "12.40"
Set-TextValue "E19" "  +5.22%  "
Set-TextValue "D20" "8.14"
Set-TextValue "E20" "  +2.75%  "
Set-TextValue "D21" "370.47"
Set-TextValue "E21" "  -3.93%  "
Set-TextValue "E22" "  +0.26%  "
Set-TextValue "E23" "  +1.65%  "
Set-TextValue "D24" "71.95"
Set-TextValue "E24" "  -2.86%  "
Set-TextValue "E25" "  +0.04%  "
Set-TextValue "D26" "4.32"
Set-TextValue "E26" "  -2.33%  "
Set-TextValue "D27" "9.82"
Set-TextValue "E27" "  -2.12%  "
Set-TextValue "D28" "2.798.78"
Set-TextValue "E28" "  +2.51%  "
Set-TextValue "D29" "1.00"
Set-TextValue "E29" "  +0.09%  "
Set-TextValue "E30" "  +0.52%  "
Set-TextValue "D31" "8.16"
Set-TextValue "E31" "  +0.48%  "
Set-TextValue "D32" "495.51"
Set-TextValue "E32" "  -4.75%  "
Set-TextValue "E33" "  -2.45%  "
Set-TextValue "D34" "1.82"
Set-TextValue "E34" "  -0.44%  "
Set-TextValue "D35" "1.00"
Set-TextValue "E35" "  +0.13%  "
Set-TextValue "E36" "  -0.11%  "
Set-TextValue "D37" "19.49"
Set-TextValue "E37" "  +0.28%  "
Set-TextValue "E38" "  +0.50%  "
Set-TextValue "E39" "  -0.93%  "
Set-TextValue "E40" "  -2.17%  "
Set-TextValue "E41" "  -5.68%  "
Set-TextValue "B43" "dogwifhat"
Set-TextValue "C43" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D43" "2.59"
Set-TextValue "E43" "  -0.62%  "
Set-TextValue "B44" "RenderToken"
Set-TextValue "C44" "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D44" "4.99"
Set-TextValue "E44" "  -3.27%  "
Set-TextValue "E45" "  -0.51%  "
Set-TextValue "D46" "156.41"
Set-TextValue "E46" "  +3.40%  "
Set-TextValue "D47" "39.29"
Set-TextValue "E47" "  -0.45%  "
Set-TextValue "E48" "  +2.32%  "
Set-TextValue "E49" "  +0.83%  "
Set-TextValue "D50" "1.72"
Set-TextValue "E50" "  +1.37%  "
Set-TextValue "B51" "BabyDogeCoin"
Set-TextValue "C51" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D51" "0.0₆0259"
Set-TextValue "E51" "  -2.45%  "
